# Update the "Förändrad" (Changed) date column (C) for rows 2-52
# from serial date 45182 (2023-09-13) to serial date 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}
